$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data rows (row 2 and row 3) for columns B through H,
# matching the updated schedule values from the source data.
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = -5
$ws.Range("H2").Value = 11

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16
